$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 3
$ws.Range("I3").Value = 2.15
$ws.Range("J3").Value = 3.5
$ws.Range("K3").Value = 2.3
$ws.Range("L3").Value = 2.75
$ws.Range("AF3").Value = 34
$ws.Range("AG3").Value = 23
$ws.Range("AO3").Value = 9

# Row 4
$ws.Range("U4").Value = 1.95
$ws.Range("V4").Value = 1.9

# Row 5
$ws.Range("H5").Value = 3.35
$ws.Range("I5").Value = 2.5
$ws.Range("J5").Value = 3.1
$ws.Range("K5").Value = 2.15
$ws.Range("L5").Value = 3
$ws.Range("O5").Value = 1.23
$ws.Range("P5").Value = 3.4
$ws.Range("S5").Value = 1.7
$ws.Range("T5").Value = 1.93
$ws.Range("W5").Value = 2.57
$ws.Range("X5").Value = 1.38
$ws.Range("AA5").Value = 1.55
$ws.Range("AB5").Value = 2.15
$ws.Range("AC5").Value = 10
$ws.Range("AD5").Value = 14
$ws.Range("AF5").Value = 29
$ws.Range("AG5").Value = 20
$ws.Range("AH5").Value = 25
$ws.Range("AI5").Value = 12
$ws.Range("AJ5").Value = 6.7
$ws.Range("AK5").Value = 12
$ws.Range("AL5").Value = 45
$ws.Range("AM5").Value = 10
$ws.Range("AS5").Value = 300

# Row 8
$ws.Range("G8").Value = 1.8
$ws.Range("I8").Value = 4
$ws.Range("J8").Value = 2.4
$ws.Range("L8").Value = 4.5
$ws.Range("M8").Value = 1.03
$ws.Range("N8").Value = 9.5
$ws.Range("O8").Value = 1.3
$ws.Range("P8").Value = 3.4
$ws.Range("S8").Value = 2
$ws.Range("T8").Value = 1.8
$ws.Range("Y8").Value = 1.4
$ws.Range("Z8").Value = 2.75
$ws.Range("AG8").Value = 15
$ws.Range("AM8").Value = 11
$ws.Range("AN8").Value = 21
$ws.Range("AO8").Value = 15
$ws.Range("AS8").Value = 800
